$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B1").Value = 11
$ws.Range("B22").Value = 28
$ws.Range("B23").Value = 980000
$ws.Range("B24").Value = 4000000
$ws.Range("B34").Value = 130000
$ws.Range("B35").Value = 130000
